# CheckLot2.xlsx - Seller sheet: change the five "LIMIT 10" lookup queries
# in column G (rows 2-6) to "LIMIT 1", and move the active selection to G2.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = "SELECT LTLOTNBR, LTSLRNBR FROM MISPRDDB.LOTLT a JOIN MISPRDDB.LOTLTNEW b ON a.LTLOTNBR = b.LOT_NUMBER WHERE LTLOTSTG =10 AND EW_TS03 IS NULL AND EW_TS04 is NULL ORDER BY LTLOTNBR DESC LIMIT 1"
$ws.Range("G3").Value = "SELECT LTLOTNBR, LTSLRNBR FROM MISPRDDB.LOTLT a JOIN MISPRDDB.LOTLTNEW b ON a.LTLOTNBR = b.LOT_NUMBER WHERE LTLOTSTG =15 AND EW_TS03 IS NULL AND EW_TS04 is NULL ORDER BY LTLOTNBR DESC LIMIT 1"
$ws.Range("G4").Value = "SELECT LTLOTNBR, LTSLRNBR FROM MISPRDDB.LOTLT a JOIN MISPRDDB.LOTLTNEW b ON a.LTLOTNBR = b.LOT_NUMBER WHERE LTLOTSTG =20 AND EW_TS03 IS NULL AND EW_TS04 is NULL ORDER BY LTLOTNBR DESC LIMIT 1"
$ws.Range("G5").Value = "SELECT LTLOTNBR, LTSLRNBR FROM MISPRDDB.LOTLT a JOIN MISPRDDB.LOTLTNEW b ON a.LTLOTNBR = b.LOT_NUMBER WHERE LTLOTSTG =28 AND EW_TS03 IS NULL AND EW_TS04 is NULL ORDER BY LTLOTNBR DESC LIMIT 1"
$ws.Range("G6").Value = "SELECT LTLOTNBR, LTSLRNBR FROM MISPRDDB.LOTLT a JOIN MISPRDDB.LOTLTNEW b ON a.LTLOTNBR = b.LOT_NUMBER WHERE LTLOTSTG =28 AND EW_TS03 IS NOT NULL AND EW_TS04 is NULL ORDER BY LTLOTNBR DESC LIMIT 1"

# Move / restore the active cell selection to G2 (was G7 before the edit).
$ws.Range("G2").Select()
